$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (invoice z0bug.invoice_PY_2 / z0bug.res_partner_2): update origin (D)
# and reference (E) from "Telefonico" to "Contratto".
$ws.Range("D3").Value = "Contratto"
$ws.Range("E3").Value = "Contratto"

# Match the author's cursor position left in the saved file.
$ws.Range("M3").Select()
